$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number ("26.00", "0.7916", ...)
# must be pre-formatted as Text so Excel keeps the literal digits/
# trailing zeros instead of silently parsing them into a Double.
$ws.Range('D2').Value = '29.858.80'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '1.892.34'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7916'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '242.39'
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('E8').Value = '  +2.63%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '26.00'
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07088'
$ws.Range('E10').Value = '  +3.08%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08064'
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7733'
$ws.Range('E12').Value = '  +5.23%  '
$ws.Range('D13').Value = '1.880.98'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.325'
$ws.Range('E14').Value = '  +3.05%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '92.37'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').Value = '29.863.72'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.927'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '244.01'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007742'
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.156.09'
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.063'
$ws.Range('E23').Value = '  +16.68%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1629'
$ws.Range('E25').Value = '  +14.56%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.334'
$ws.Range('E26').Value = '  +1.79%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.13'
$ws.Range('E27').Value = '  -1.32%  '
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('E29').Value = '  +1.70%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.383'
$ws.Range('E30').Value = '  +2.18%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.536'
$ws.Range('E31').Value = '  +1.60%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.433'
$ws.Range('E32').Value = '  +3.48%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05644'
$ws.Range('E33').Value = '  +1.19%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.103'
$ws.Range('E34').Value = '  +1.18%  '
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7372'
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.003'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.709'
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01932'
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.776'
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4458'
$ws.Range('E41').Value = '  +1.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '72.29'
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.871'
$ws.Range('E43').Value = '  -2.03%  '
$ws.Range('E44').Value = '  +1.43%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.002'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.888'
$ws.Range('E46').Value = '  +1.57%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '102.46'
$ws.Range('E47').Value = '  +2.04%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.023.10'
$ws.Range('E48').Value = '  +4.79%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.889'
$ws.Range('E49').Value = '  +1.88%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.494'
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.969'
$ws.Range('E51').Value = '  +7.07%  '
